# Generate Report for Handback
#
# Populates row 7 (the "9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9" source file)
# on both localized-status sheets (zh-cn / de-de) now that a handback has
# been processed for it: the "Latest Target File", "Latest Handback File",
# "Latest Handback DateTime" and "Error Detail" columns move from their
# placeholder/empty state to the real handback report values.

$wb = $excel.ActiveWorkbook

$latestUrl = "https://github.com/OpenLocalizationTestOrg/oltest/blob/689a96aa9dc25861323a9297e2bc1d4a69b0bac9/e2e/9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.md"
$targetDisplay = "9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.md"
$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/bbc923119251aa304fcaf6ced870f1da36b491ee/e2e/9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/689a96aa9dc25861323a9297e2bc1d4a69b0bac9/e2e/9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.md."

# ---------------------------------------------------------------------
# zh-cn sheet, row 7
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$iCell = $wsZhCn.Range("I7")
$iCell.Value = $targetDisplay
$wsZhCn.Hyperlinks.Add($iCell, $latestUrl, [System.Type]::Missing, [System.Type]::Missing, $targetDisplay) | Out-Null

$wsZhCn.Range("J7").Value = "9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.9243079b6d47087028949c9654da031563f83376.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-12 21:07:18"
$wsZhCn.Range("P7").Value = $errorMessage

# ---------------------------------------------------------------------
# de-de sheet, row 7
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$iCell2 = $wsDeDe.Range("I7")
$iCell2.Value = $targetDisplay
$wsDeDe.Hyperlinks.Add($iCell2, $latestUrl, [System.Type]::Missing, [System.Type]::Missing, $targetDisplay) | Out-Null

$wsDeDe.Range("J7").Value = "9fb5e8a5-9865-4df0-8d2b-9a90de81e3e9.9243079b6d47087028949c9654da031563f83376.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-12 21:07:27"
$wsDeDe.Range("P7").Value = $errorMessage
